$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.513.13"
$ws.Range("E2").Value = "  -1.62%  "
$ws.Range("D3").Value = "1.628.47"
$ws.Range("E3").Value = "  -0.76%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.21"
$ws.Range("E5").Value = "  -0.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.521"
$ws.Range("E6").Value = "  -1.03%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.91"
$ws.Range("E8").Value = "  -1.85%  "
$ws.Range("E9").Value = "  -0.57%  "
$ws.Range("E10").Value = "  -0.35%  "
$ws.Range("E11").Value = "  -3.35%  "
$ws.Range("D12").Value = "1.860.03"
$ws.Range("E12").Value = "  -0.70%  "
$ws.Range("D13").Value = "1.624.51"
$ws.Range("E13").Value = "  -0.97%  "
$ws.Range("E14").Value = "  -0.39%  "
$ws.Range("E15").Value = "  -0.59%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.05"
$ws.Range("E16").Value = "  +0.57%  "
$ws.Range("D17").Value = "27.513.64"
$ws.Range("E17").Value = "  -1.51%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "228.59"
$ws.Range("E18").Value = "  -2.06%  "
$ws.Range("D19").Value = "0.0₃0717"
$ws.Range("E19").Value = "  -0.75%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.52"
$ws.Range("E20").Value = "  -1.47%  "
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.73"
$ws.Range("E22").Value = "  +7.37%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.37"
$ws.Range("E23").Value = "  +1.26%  "
$ws.Range("E24").Value = "  +2.87%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.86"
$ws.Range("E25").Value = "  -1.25%  "
$ws.Range("E26").Value = "  -1.25%  "
$ws.Range("E27").Value = "  -0.69%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.58"
$ws.Range("E28").Value = "  -0.82%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("E30").Value = "  -0.82%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0480"
$ws.Range("E31").Value = "  -0.58%  "
$ws.Range("E32").Value = "  -1.39%  "
$ws.Range("D33").Value = "1.461.91"
$ws.Range("E33").Value = "  -0.73%  "
$ws.Range("E34").Value = "  -0.69%  "
$ws.Range("E35").Value = "  -0.68%  "
$ws.Range("E36").Value = "  -1.63%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.875"
$ws.Range("E37").Value = "  -0.54%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.557"
$ws.Range("E38").Value = "  -1.64%  "
$ws.Range("E39").Value = "  -0.69%  "
$ws.Range("E40").Value = "  -0.06%  "
$ws.Range("B41").Value = "WEMIXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.02"
$ws.Range("E41").Value = "  +0.64%  "
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "68.19"
$ws.Range("E43").Value = "  -1.48%  "
$ws.Range("E44").Value = "  +0.25%  "
$ws.Range("E45").Value = "  -1.15%  "
$ws.Range("E46").Value = "  -1.24%  "
$ws.Range("D47").Value = "1.768.92"
$ws.Range("E47").Value = "  -0.77%  "
$ws.Range("E48").Value = "  +1.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "87.32"
$ws.Range("E49").Value = "  +0.36%  "
$ws.Range("E50").Value = "  -1.24%  "
$ws.Range("E51").Value = "  +0.17%  "
